$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill in the "VIERNES" column (G) for the second member's stand-up rows (7-9),
# matching the text already used for "Nada" / the new weekly-meeting note.
$ws.Range("G7").Value = "Nada"
$ws.Range("G8").Value = "Reunion semanal y asignacion de tareas"
$ws.Range("G9").Value = "Nada"

# Widen column G so the new text is readable.
$ws.Columns.Item(7).ColumnWidth = 37.75

# Update the active sheet view: scrolled one column further right (topLeftCell
# D1 -> E1), and the current selection moved from J10 to G9.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("G9").Select()
